{"js": "// \"Version 2.\" -> \"Version 1.\"\n// (restores the pre-\"Wireframes version 2\" wording, per the commit message:\n// \"Revert \"Revert \"Revert \"Wireframes version 2.\"\"\"\")\nconst body = context.document.body;\n\n// 1) Re-type \"Version\" over itself so the pre-existing \"Versi\"/\"on\" run\n//    split (left over from an earlier spell-check/autocorrect edit)\n//    collapses into a single run, matching how Word naturally coalesces\n//    runs when a span is retyped.\nconst versionMatches = body.search(\"Version\", { matchCase: true, matchWholeWord: false });\nversionMatches.load(\"text\");\nawait context.sync();\nif (versionMatches.items.length > 0) {\n  versionMatches.items[0].insertText(\"Version\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 2) Replace \" 2.\" (space + digit + trailing period) with \" 1.\" \u2014 this\n//    both updates the version number and absorbs the separate trailing\n//    \".\" run into the number's run, just like the target markup.\nconst numberMatches = body.search(\" 2.\", { matchCase: true, matchWholeWord: false });\nnumberMatches.load(\"text\");\nawait context.sync();\nif (numberMatches.items.length > 0) {\n  numberMatches.items[0].insertText(\" 1.\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# \"Version 2.\" -> \"Version 1.\"\n# (restores the pre-\"Wireframes version 2\" wording, per the commit message:\n# Revert \"Revert \"Revert \"Wireframes version 2.\"\"\")\n$d = $word.ActiveDocument\n\n# 1) Re-type \"Version\" over itself so the pre-existing \"Versi\"/\"on\" run\n#    split (left over from an earlier spell-check/autocorrect edit)\n#    collapses into a single run, matching how Word naturally coalesces\n#    runs when a span is retyped.\n$find0 = $d.Content.Find\n$find0.ClearFormatting()\n$find0.Execute(\"Version\", $false, $true, $false, $false, $false, $true, 1, $false, \"Version\", 2)\n\n# 2) Replace the digit \"2\" with \"1.\" (does not cross the \"_GoBack\" bookmark\n#    that sits right after the digit, so the bookmark survives). This both\n#    updates the version number and gives this run its own trailing period.\n$find1 = $d.Range().Find\n$find1.ClearFormatting()\n$find1.Execute(\"2\", $false, $false, $false, $false, $false, $true, 1, $false, \"1.\", 2)\n\n# 3) The original trailing \".\" run (now duplicated, right after the\n#    bookmark) is deleted by character position rather than Find, again so\n#    the bookmark is left untouched. Content.End points just past the final\n#    paragraph mark, so Content.End - 2 is the last real character (\".\").\n$docEnd = $d.Content.End\n$d.Range($docEnd - 2, $docEnd - 1).Delete()\n"}
